# DaySale report update:
#  - item 1 (row 7) changes from "صابون ديتول اوريجنيال 115 جم" to "DECONGESTYL-N 12 RECTAL SUPP."
#  - two new item rows are inserted after it (LARYPRO 20 LOZENGES, then the original
#    "صابون ديتول اوريجنيال 115 جم" row moves down to become item 3)
#  - the totals row and the footer (generated-at timestamp) row move down accordingly
#    and the timestamp is refreshed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the old totals row (row 8), shifting the totals/footer
# rows down to 10/11 and carrying their merges/styles with them.
$ws.Rows("8:9").Insert()

# Row 8 and row 9 need the same per-column formatting as row 7 (the item-row template).
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)
$ws.Range("A7:Q7").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)
$ws.Range("Q7").Copy()
$ws.Range("Q8").PasteSpecial(-4122)
$ws.Range("Q7").Copy()
$ws.Range("Q9").PasteSpecial(-4122)

# Recreate the merges row 7 has, for the two new item rows.
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

# Row heights as in the target layout.
$ws.Rows("8").RowHeight = 24.75
$ws.Rows("9").RowHeight = 25.5

function Set-TextValue($rng, $text) {
    # Write a text value into a cell whose number format is numeric, without
    # Excel auto-converting the literal to a real number (and without leaving
    # the "number stored as text" quote-prefix marker behind).
    $fmt = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = $fmt
}

# --- Item 1 (row 7): now DECONGESTYL-N 12 RECTAL SUPP. ---
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "DECONGESTYL-N 12 RECTAL SUPP."
$ws.Range("H7").Value = "0:1"
Set-TextValue $ws.Range("L7") "1"
$ws.Range("N7").Value = "39.00"
Set-TextValue $ws.Range("P7") "39.0000"
$ws.Range("Q7").Value = "1:0"

# --- Item 2 (row 8): LARYPRO 20 LOZENGES ---
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "LARYPRO 20 LOZENGES"
$ws.Range("H8").Value = "0:1"
Set-TextValue $ws.Range("L8") "1"
$ws.Range("N8").Value = "44.00"
Set-TextValue $ws.Range("P8") "22.0000"
$ws.Range("Q8").Value = "0:1"

# --- Item 3 (row 9): صابون ديتول اوريجنيال 115 جم (the original item 1) ---
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "صابون ديتول اوريجنيال 115 جم"
$ws.Range("H9").Value = "0:0"
Set-TextValue $ws.Range("L9") "0"
$ws.Range("N9").Value = "30.00"
Set-TextValue $ws.Range("P9") "30.0000"
$ws.Range("Q9").Value = "1:0"

# --- Totals row (now row 10) ---
$ws.Range("P10").Value = 91

# --- Footer row (now row 11): refreshed generation timestamp ---
$ws.Range("A11").Value = "Monday, 28 July, 2025 9:43 AM"

$wb.Save()
